$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3562
$ws.Range("F5").Value = 8328
$ws.Range("F7").Value = 124
$ws.Range("F8").Value = 2229
$ws.Range("F9").Value = 18
$ws.Range("F10").Value = 101
$ws.Range("F11").Value = 74
$ws.Range("F12").Value = 649
$ws.Range("F13").Value = 105
$ws.Range("F14").Value = 7449
$ws.Range("F16").Value = 7691
$ws.Range("F18").Value = 57755
$ws.Range("F19").Value = 57756
$ws.Range("F20").Value = 4796
$ws.Range("F21").Value = 1058
$ws.Range("F22").Value = 947
$ws.Range("F23").Value = 502
$ws.Range("F24").Value = 111
$ws.Range("F28").Value = 5306
$ws.Range("F30").Value = 118
$ws.Range("F32").Value = 912
$ws.Range("F33").Value = 1392
$ws.Range("F34").Value = 1948
$ws.Range("F36").Value = 185
$ws.Range("F38").Value = 1088
$ws.Range("F42").Value = 783
$ws.Range("F43").Value = 270
$ws.Range("F44").Value = 209
$ws.Range("F47").Value = 205
$ws.Range("F49").Value = 59
$ws.Range("F50").Value = 2487

# ---- Sheet: 演出 ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 36
$ws.Range("F6").Value = 140
$ws.Range("F8").Value = 50
$ws.Range("F9").Value = 7633
$ws.Range("F10").Value = 126
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 4
$ws.Range("F19").Value = 22
$ws.Range("F22").Value = 33
$ws.Range("F24").Value = 3
$ws.Range("F31").Value = 81
$ws.Range("F32").Value = 1
$ws.Range("F45").Value = 29
$ws.Range("F48").Value = 279

# ---- Sheet: 本地生活 ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2373
$ws.Range("F5").Value = 1608
$ws.Range("E7").Value = "2024.08.16 00:00-10.13 23:59"
$ws.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202409/hAWFKrvi1727676771939.png"
$ws.Range("F8").Value = 2422
$ws.Range("F9").Value = 9439
$ws.Range("F10").Value = 1769
$ws.Range("F11").Value = 181
$ws.Range("F15").Value = 271
$ws.Range("F16").Value = 2375
$ws.Range("C17").Value = "上海·三丽鸥家族Sanrio Characters主题餐厅·海滩奇遇季"
$ws.Range("F17").Value = 29
$ws.Range("F18").Value = 64
$ws.Range("F19").Value = 518

# ---- Sheet: 全部类型 ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3562
$ws.Range("F3").Value = 2373
$ws.Range("F5").Value = 8328
$ws.Range("E6").Value = "2024.08.16 00:00-10.13 23:59"
$ws.Range("I6").Value = "//i0.hdslb.com/bfs/openplatform/202409/hAWFKrvi1727676771939.png"
$ws.Range("F7").Value = 181
$ws.Range("F8").Value = 124
$ws.Range("F9").Value = 271
$ws.Range("F10").Value = 74
$ws.Range("F11").Value = 649
$ws.Range("F12").Value = 105
$ws.Range("F13").Value = 7449
$ws.Range("F14").Value = 7692
$ws.Range("F15").Value = 57756
$ws.Range("F17").Value = 36
$ws.Range("F18").Value = 4796
$ws.Range("F19").Value = 1058
$ws.Range("F20").Value = 947
$ws.Range("F21").Value = 502
$ws.Range("F23").Value = 140
$ws.Range("F24").Value = 5306
$ws.Range("F26").Value = 118
$ws.Range("F27").Value = 912
$ws.Range("F28").Value = 1392
$ws.Range("F29").Value = 126
$ws.Range("F30").Value = 518
$ws.Range("F31").Value = 4
$ws.Range("F33").Value = 185
$ws.Range("F34").Value = 22
$ws.Range("F36").Value = 783
$ws.Range("F37").Value = 270
$ws.Range("F40").Value = 3
$ws.Range("F45").Value = 205
$ws.Range("F49").Value = 2487
$ws.Range("F50").Value = 29
